$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 0
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -4
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = 6
$ws.Range("F13").Value = -2
$ws.Range("F18").Value = -3
